$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.961.62'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -3.20%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.860.73'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -2.44%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '317.97'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.07%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -4.69%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3706'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -2.86%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07512'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -2.72%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9384'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -4.25%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '21.28'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -4.01%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.878.09'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.733'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -3.15%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.446'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -4.09%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06831'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -3.23%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '81.71'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -2.44%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000009023'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -4.71%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.9998'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.96'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -4.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '27.950.01'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -3.24%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.115'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -3.87%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.07'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.091.79'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.008'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -4.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '154.65'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -2.60%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.42'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -3.08%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.438'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -4.08%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '113.61'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -3.53%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.743'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -6.74%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08996'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -3.28%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.8137'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -5.84%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.825'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -5.18%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.176'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -5.61%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.932'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.75%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.9999'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05499'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -3.74%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.114'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.66%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01978'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -3.21%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.903'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.5264'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -4.13%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '7.050'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -5.40%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1692'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -3.49%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.800'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -5.80%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.06785'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4910'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -4.90%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.58'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -5.88%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '106.30'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.682'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -5.33%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.911'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -11.75%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9991'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.23%  '
